$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp label (row 1)
$ws.Range("A1").Value = "Datos actualizados a 17 de Junio de 2020 a las 18:49"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 2217920
$ws.Range("C4").Value = 9520
$ws.Range("D4").Value = 903616
$ws.Range("E4").Value = 1194949
$ws.Range("G4").Value = 223
$ws.Range("H4").Value = 119355

# --- Row 5: Brasil ---
$ws.Range("B5").Value = 934769
$ws.Range("C5").Value = 5935
$ws.Range("E5").Value = 411820
$ws.Range("G5").Value = 129
$ws.Range("H5").Value = 45585

# --- Row 7: India ---
$ws.Range("B7").Value = 360365
$ws.Range("C7").Value = 6204
$ws.Range("D7").Value = 191257
$ws.Range("E7").Value = 157051
$ws.Range("G7").Value = 136
$ws.Range("H7").Value = 12057

# --- Row 10: Italia ---
$ws.Range("B10").Value = 237828
$ws.Range("C10").Value = 328
$ws.Range("D10").Value = 179455
$ws.Range("E10").Value = 23925
$ws.Range("G10").Value = 43
$ws.Range("H10").Value = 34448

# --- Row 34: Singapur ---
$ws.Range("D34").Value = 31938
$ws.Range("E34").Value = 9252

# --- Row 56: Kazajistan ---
$ws.Range("D56").Value = 9920
$ws.Range("E56").Value = 5534

# --- Rows 57-60: countries re-ranked (Ghana/Serbia/Dinamarca/Moldavia shift) ---
$ws.Range("A57").Value = "Moldavia"
$ws.Range("B57").Value = 12732
$ws.Range("C57").Value = 478
$ws.Range("D57").Value = 7077
$ws.Range("E57").Value = 5222
$ws.Range("G57").Value = 10
$ws.Range("H57").Value = 433

$ws.Range("A58").Value = "Ghana"
$ws.Range("B58").Value = 12590
$ws.Range("C58").Value = 397
$ws.Range("D58").Value = 4410
$ws.Range("E58").Value = 8114
$ws.Range("G58").Value = 8
$ws.Range("H58").Value = 66

$ws.Range("A59").Value = "Serbia"
$ws.Range("B59").Value = 12522
$ws.Range("C59").Value = 96
$ws.Range("D59").Value = 11511
$ws.Range("E59").Value = 754
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 257

$ws.Range("A60").Value = "Dinamarca"
$ws.Range("B60").Value = 12294
$ws.Range("C60").Value = 44
$ws.Range("D60").Value = 11185
$ws.Range("E60").Value = 511
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 598

# --- Row 62: Argelia ---
$ws.Range("B62").Value = 11268
$ws.Range("C62").Value = 121
$ws.Range("D62").Value = 7943
$ws.Range("E62").Value = 2526
$ws.Range("G62").Value = 11
$ws.Range("H62").Value = 799

# --- Row 65: Chequia ---
$ws.Range("B65").Value = 10154
$ws.Range("C65").Value = 43
$ws.Range("D65").Value = 7399
$ws.Range("E65").Value = 2422
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 333

# --- Row 82: Republica de Yibuti ---
$ws.Range("B82").Value = 4545
$ws.Range("C82").Value = 6
$ws.Range("D82").Value = 3411
$ws.Range("E82").Value = 1091

# --- Rows 85-86: Hungria/Luxemburgo swap ---
$ws.Range("A85").Value = "Luxemburgo"
$ws.Range("B85").Value = 4085
$ws.Range("C85").Value = 10
$ws.Range("D85").Value = 3935
$ws.Range("E85").Value = 40
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 110

$ws.Range("A86").Value = "Hungria"
$ws.Range("B86").Value = 4078
$ws.Range("C86").Value = 1
$ws.Range("D86").Value = 2547
$ws.Range("E86").Value = 964
$ws.Range("G86").Value = 2
$ws.Range("H86").Value = 567

# --- Row 88: Kenia ---
$ws.Range("B88").Value = 4044
$ws.Range("C88").Value = 184
$ws.Range("D88").Value = 1353
$ws.Range("E88").Value = 2584
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 107

# --- Rows 206-207: Groenlandia/Islas Malvinas swap (values identical, only names swap) ---
$ws.Range("A206").Value = "Islas Malvinas"
$ws.Range("A207").Value = "Groenlandia"

# --- Rows 210-211: Seychelles/Montserrat swap ---
$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# --- Rows 213-214: Papua Nueva Guinea/Islas Virgenes Britanicas swap ---
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
